$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.405.87"
$ws.Range("E2").Value = "  -2.85%  "
$ws.Range("D3").Value = "1.984.62"
$ws.Range("E3").Value = "  -3.53%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.53"
$ws.Range("E5").Value = "  -3.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.630"
$ws.Range("E6").Value = "  -3.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.84"
$ws.Range("E7").Value = "  -11.67%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  -1.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.29"
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("E11").Value = "  +7.57%  "
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.69"
$ws.Range("E13").Value = "  +5.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.869"
$ws.Range("E14").Value = "  -6.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.05"
$ws.Range("E15").Value = "  -5.59%  "
$ws.Range("D16").Value = "2.273.50"
$ws.Range("E16").Value = "  -3.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.49"
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").Value = "1.980.00"
$ws.Range("E18").Value = "  -3.62%  "
$ws.Range("D19").Value = "36.285.15"
$ws.Range("E19").Value = "  -2.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.39"
$ws.Range("E20").Value = "  -4.48%  "
$ws.Range("D21").Value = "0.0₃0877"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.34"
$ws.Range("E22").Value = "  -3.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.58"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  -5.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.31"
$ws.Range("E26").Value = "  -4.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.97"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.92"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.86"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.131"
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.92"
$ws.Range("E33").Value = "  -6.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0651"
$ws.Range("E34").Value = "  +3.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.43"
$ws.Range("E35").Value = "  -5.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.24"
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.27"
$ws.Range("E38").Value = "  -7.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.78"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.03"
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.24"
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0973"
$ws.Range("E42").Value = "  -6.25%  "
$ws.Range("E43").Value = "  -4.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0214"
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("E45").Value = "  -5.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.29"
$ws.Range("E46").Value = "  -7.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "92.59"
$ws.Range("E47").Value = "  -5.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.51"
$ws.Range("E48").Value = "  -5.53%  "
$ws.Range("D49").Value = "1.370.73"
$ws.Range("E49").Value = "  -3.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.83"
$ws.Range("E50").Value = "  -4.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.66"
$ws.Range("E51").Value = "  -2.25%  "
